{"js": "// The document contains a single table of 20 rows x 5 columns (100\n// cells total), each holding one \"a op b = c\" arithmetic-practice\n// answer. The edit rewrites all 100 answers in place (row-major order:\n// row 0..19, left-to-right column 0..4), matching the canonical-XML\n// diff which swaps each cell's <w:t> run text one-for-one without\n// touching any other markup/formatting.\nconst newValues = [\n  \"41+46=87\", \"11+71=82\", \"36+37=73\", \"61-47=14\", \"35+1=36\",\n  \"94-76=18\", \"67-30=37\", \"97-92=5\", \"80+8=88\", \"6+6=12\",\n  \"65-19=46\", \"75+4=79\", \"70+12=82\", \"65-16=49\", \"66+21=87\",\n  \"3+75=78\", \"48-16=32\", \"63-54=9\", \"54-35=19\", \"25-11=14\",\n  \"28-17=11\", \"96-67=29\", \"63-43=20\", \"16+0=16\", \"20+21=41\",\n  \"95-47=48\", \"26+41=67\", \"90-0=90\", \"66+3=69\", \"21+70=91\",\n  \"59-27=32\", \"95+3=98\", \"53-44=9\", \"37+11=48\", \"24+40=64\",\n  \"15+61=76\", \"81-42=39\", \"2+52=54\", \"68-39=29\", \"93-18=75\",\n  \"60-23=37\", \"37+61=98\", \"26+59=85\", \"7+8=15\", \"13+43=56\",\n  \"54+43=97\", \"12+15=27\", \"8+91=99\", \"7+88=95\", \"90-37=53\",\n  \"59-7=52\", \"0+77=77\", \"70-61=9\", \"67-16=51\", \"2+85=87\",\n  \"62-41=21\", \"40-6=34\", \"2+11=13\", \"34+16=50\", \"35+29=64\",\n  \"16-13=3\", \"13+52=65\", \"19+2=21\", \"50-2=48\", \"46-30=16\",\n  \"30+66=96\", \"47+11=58\", \"91-55=36\", \"49+10=59\", \"45+54=99\",\n  \"73-19=54\", \"20+21=41\", \"45+45=90\", \"1+28=29\", \"41+47=88\",\n  \"81-56=25\", \"54-14=40\", \"26+46=72\", \"69+7=76\", \"19+20=39\",\n  \"53-37=16\", \"50-26=24\", \"70-45=25\", \"31-15=16\", \"50+39=89\",\n  \"76-20=56\", \"72-27=45\", \"88-83=5\", \"26+58=84\", \"54+44=98\",\n  \"98-78=20\", \"19+80=99\", \"28+53=81\", \"29+29=58\", \"90-14=76\",\n  \"83-32=51\", \"10+81=91\", \"85+2=87\", \"4+32=36\", \"70+4=74\",\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 5;\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (i >= newValues.length) break;\n    table.getCell(r, c).value = newValues[i];\n    i++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-equation answers in the single table,\n# cell by cell, in row-major order (row 1..20, column 1..5). This\n# mirrors the canonical-XML diff, which rewrites each <w:t> run's text\n# content one-for-one, in document order, without touching any other\n# formatting/markup.\n$values = @(\n  \"41+46=87\", \"11+71=82\", \"36+37=73\", \"61-47=14\", \"35+1=36\", \"94-76=18\", \"67-30=37\", \"97-92=5\", \"80+8=88\", \"6+6=12\", \"65-19=46\", \"75+4=79\", \"70+12=82\", \"65-16=49\", \"66+21=87\", \"3+75=78\", \"48-16=32\", \"63-54=9\", \"54-35=19\", \"25-11=14\", \"28-17=11\", \"96-67=29\", \"63-43=20\", \"16+0=16\", \"20+21=41\", \"95-47=48\", \"26+41=67\", \"90-0=90\", \"66+3=69\", \"21+70=91\", \"59-27=32\", \"95+3=98\", \"53-44=9\", \"37+11=48\", \"24+40=64\", \"15+61=76\", \"81-42=39\", \"2+52=54\", \"68-39=29\", \"93-18=75\", \"60-23=37\", \"37+61=98\", \"26+59=85\", \"7+8=15\", \"13+43=56\", \"54+43=97\", \"12+15=27\", \"8+91=99\", \"7+88=95\", \"90-37=53\", \"59-7=52\", \"0+77=77\", \"70-61=9\", \"67-16=51\", \"2+85=87\", \"62-41=21\", \"40-6=34\", \"2+11=13\", \"34+16=50\", \"35+29=64\", \"16-13=3\", \"13+52=65\", \"19+2=21\", \"50-2=48\", \"46-30=16\", \"30+66=96\", \"47+11=58\", \"91-55=36\", \"49+10=59\", \"45+54=99\", \"73-19=54\", \"20+21=41\", \"45+45=90\", \"1+28=29\", \"41+47=88\", \"81-56=25\", \"54-14=40\", \"26+46=72\", \"69+7=76\", \"19+20=39\", \"53-37=16\", \"50-26=24\", \"70-45=25\", \"31-15=16\", \"50+39=89\", \"76-20=56\", \"72-27=45\", \"88-83=5\", \"26+58=84\", \"54+44=98\", \"98-78=20\", \"19+80=99\", \"28+53=81\", \"29+29=58\", \"90-14=76\", \"83-32=51\", \"10+81=91\", \"85+2=87\", \"4+32=36\", \"70+4=74\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($i -ge $values.Count) { break }\n    $t.Cell($r, $c).Range.Text = $values[$i]\n    $i++\n  }\n}\n"}
